# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / handoff / handback timestamps
# for the 974ea617-c12b-41f4-a54d-b6162e8832ff.md entry (row 3 in every sheet)
# to reflect a fresh handback report generation.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-11-03 19:27:04"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-11-03 19:26:51"
$zhcn.Range("K3").Value = "2016-11-03 19:27:42"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-11-03 19:27:04"
$dede.Range("K3").Value = "2016-11-03 19:27:59"
